$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 16667010
$ws.Range("I33").Value = 25000240
$ws.Range("J33").Value = 550
$ws.Range("K33").Value = 25000240
$ws.Range("L33").Value = 550
$ws.Range("M33").Value = -25000011
$ws.Range("N33").Value = -1008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2054.75
$ws.Range("I40").Value = 1379.5
$ws.Range("J40").Value = 2279.8333
$ws.Range("K40").Value = 1379.5
$ws.Range("L40").Value = 2279.8333
$ws.Range("M40").Value = -1204.5
$ws.Range("N40").Value = -2629.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 18570.166
$ws.Range("I43").Value = 34500
$ws.Range("K43").Value = 34500
$ws.Range("M43").Value = -34431

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10619.158
$ws.Range("I116").Value = 11747.083
$ws.Range("K116").Value = 11747.083
$ws.Range("M116").Value = -8305.083000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 14136.75
$ws.Range("J137").Value = 34733.332
$ws.Range("L137").Value = 104199.996
$ws.Range("N137").Value = -109299.996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3036.476
$ws.Range("J138").Value = 3698.0217
$ws.Range("L138").Value = 11094.0651
$ws.Range("N138").Value = -21374.0651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3151.5
$ws.Range("I63").Value = 555
$ws.Range("K63").Value = 555
$ws.Range("M63").Value = 131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3151.5
$ws.Range("I66").Value = 555
$ws.Range("K66").Value = 2775
$ws.Range("M66").Value = 657

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 61444.8
$ws.Range("I74").Value = 102250.9
$ws.Range("J74").Value = 20638.7
$ws.Range("K74").Value = 102250.9
$ws.Range("L74").Value = 20638.7
$ws.Range("M74").Value = -101376.9
$ws.Range("N74").Value = -22386.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 61444.8
$ws.Range("I77").Value = 102250.9
$ws.Range("J77").Value = 20638.7
$ws.Range("K77").Value = 511254.5
$ws.Range("L77").Value = 103193.5
$ws.Range("M77").Value = -506886.5
$ws.Range("N77").Value = -111929.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 54374.75
$ws.Range("I140").Value = 47500
$ws.Range("J140").Value = 56666.332
$ws.Range("K140").Value = 47500
$ws.Range("L140").Value = 56666.332
$ws.Range("M140").Value = -42320
$ws.Range("N140").Value = -67026.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1592.1052
$ws.Range("I17").Value = 1515.3846
$ws.Range("J17").Value = 1758.3334
$ws.Range("K17").Value = 1515.3846
$ws.Range("L17").Value = 1758.3334
$ws.Range("M17").Value = -1341.3846
$ws.Range("N17").Value = -2106.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("K25").Value = 200
$ws.Range("M25").Value = -26

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21913.674
$ws.Range("I31").Value = 28442.947
$ws.Range("K31").Value = 28442.947
$ws.Range("M31").Value = -28147.947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 21913.674
$ws.Range("I34").Value = 28442.947
$ws.Range("K34").Value = 28442.947
$ws.Range("M34").Value = -28240.947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 49995
$ws.Range("J41").Value = 49995
$ws.Range("L41").Value = 49995
$ws.Range("M41").Value = -50851

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 32000
$ws.Range("J50").Value = 32000
$ws.Range("L50").Value = 32000
$ws.Range("N50").Value = -33250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 29999.5
$ws.Range("I51").Value = 29999
$ws.Range("K51").Value = 29999
$ws.Range("M51").Value = -29263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 29999.5
$ws.Range("I61").Value = 29999
$ws.Range("K61").Value = 29999
$ws.Range("M61").Value = -29651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 75072.5
$ws.Range("I127").Value = 62248.75
$ws.Range("K127").Value = 62248.75
$ws.Range("M127").Value = -57288.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 17150
$ws.Range("J133").Value = 17150
$ws.Range("L133").Value = 17150
$ws.Range("N133").Value = -22210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 997.2857
$ws.Range("I48").Value = 607.6667
$ws.Range("J48").Value = 1289.5
$ws.Range("K48").Value = 1823.0001
$ws.Range("L48").Value = 3868.5
$ws.Range("M48").Value = -1573.0001
$ws.Range("N48").Value = -4368.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 600
$ws.Range("J121").Value = 600
$ws.Range("L121").Value = 1800
$ws.Range("N121").Value = -4420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 23034.541
$ws.Range("J131").Value = 2856.25
$ws.Range("L131").Value = 8568.75
$ws.Range("N131").Value = -18648.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2370.524
$ws.Range("I137").Value = 1400.2307
$ws.Range("J137").Value = 3947.25
$ws.Range("K137").Value = 4200.6921
$ws.Range("L137").Value = 11841.75
$ws.Range("M137").Value = 899.3078999999998
$ws.Range("N137").Value = -22041.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 249.63158
$ws.Range("I2").Value = 180.375
$ws.Range("J2").Value = 368.35715
$ws.Range("K2").Value = 180.375
$ws.Range("L2").Value = 368.35715
$ws.Range("M2").Value = -67.375
$ws.Range("N2").Value = -594.35715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 3022000
$ws.Range("J21").Value = 3515000
$ws.Range("L21").Value = 3515000
$ws.Range("N21").Value = -3515346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 3022000
$ws.Range("J30").Value = 3515000
$ws.Range("L30").Value = 3515000
$ws.Range("N30").Value = -3515210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 26891.666
$ws.Range("J134").Value = 26891.666
$ws.Range("L134").Value = 80674.99800000001
$ws.Range("N134").Value = -85744.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3713.6667
$ws.Range("I7").Value = 3098.7727
$ws.Range("K7").Value = 3098.7727
$ws.Range("M7").Value = -2986.7727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 3988.389
$ws.Range("I31").Value = 232.77777
$ws.Range("J31").Value = 7744
$ws.Range("K31").Value = 232.77777
$ws.Range("L31").Value = 7744
$ws.Range("M31").Value = 15.22223
$ws.Range("N31").Value = -8240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1119.9
$ws.Range("I46").Value = 1125
$ws.Range("J46").Value = 1099.5
$ws.Range("K46").Value = 1125
$ws.Range("L46").Value = 1099.5
$ws.Range("M46").Value = -937
$ws.Range("N46").Value = -1475.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 569
$ws.Range("J55").Value = 662.5
$ws.Range("L55").Value = 662.5
$ws.Range("N55").Value = -1008.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3713.6667
$ws.Range("I126").Value = 3098.7727
$ws.Range("K126").Value = 9296.3181
$ws.Range("M126").Value = -6826.3181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2258.8572
$ws.Range("I122").Value = 1919.1052
$ws.Range("K122").Value = 5757.3156
$ws.Range("M122").Value = -3307.3156

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1857.6897
$ws.Range("I136").Value = 1638.5
$ws.Range("J136").Value = 2909.8
$ws.Range("K136").Value = 4915.5
$ws.Range("L136").Value = 8729.400000000001
$ws.Range("M136").Value = -2365.5
$ws.Range("N136").Value = -13829.4
